$p = $ppt.ActivePresentation
$m = $p.SlideMaster

# Title placeholder: "Click to edit Master title style"
$titleShape = $m.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Lines(1).Text = "ickclay otay edithay astermay itletay estylay"

# Body placeholder: "Click to edit Master text styles" + level lines
$bodyShape = $m.Shapes.Item(2)
$bodyTr = $bodyShape.TextFrame.TextRange
$bodyTr.Lines(1).Text = "ickclay otay edithay astermay exttay esstylay"
$bodyTr.Lines(2).Text = "econdsay evellay"
$bodyTr.Lines(3).Text = "irdthay evellay"
$bodyTr.Lines(4).Text = "ourthfay evellay"
$bodyTr.Lines(5).Text = "ifthfay evellay"
